$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Our Work Data" label (I12) to the new model label ---
$ws.Range("I12").Value = "Our Model eq 16 & 17"

# --- Update the critical-points table header row (row 13) ---
# New "Salts" column header + relabeled Molality/Temperature columns.
# Values are entered in an order that reproduces the shared-string table
# layout of the committed workbook (Temperature then Molality, Salts last).
$ws.Range("J13").Value = "Critical Temperature"
$ws.Range("L13").Value = "Critical Temperature"
$ws.Range("I13").Value = "Critical Molality"
$ws.Range("K13").Value = "Critical Molality"
$ws.Range("H13").Value = "Salts"

# --- Fix the Molality Deviation value for CaCl2(rich) (N14) ---
$ws.Range("N14").Value = -10.932201008676131

# --- Normalize M5's style (drop the stray font/fill/border flags) ---
$ws.Range("M5").HorizontalAlignment = -4131

# --- Update the active selection left by the author ---
$ws.Range("H12:O18").Select() | Out-Null
